$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Price values that would otherwise be auto-parsed as numbers by Excel are
# written with a leading apostrophe so they stay plain text (matching the
# original inline-string / General-format cells), same as the multi-dot
# "thousands.decimal" style prices which already stay text on their own.

$ws.Range("D2").Value = "69.498.89"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "3.758.46"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'615.40"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").Value = "'177.62"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("D7").Value = "3.757.81"
$ws.Range("E7").Value = "  -0.69%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = "  -1.57%  "
$ws.Range("D10").Value = "'0.167"
$ws.Range("E10").Value = "  -2.79%  "
$ws.Range("D11").Value = "'6.65"
$ws.Range("E11").Value = "  +4.81%  "
$ws.Range("E12").Value = "  -1.97%  "
$ws.Range("D13").Value = "'40.12"
$ws.Range("E13").Value = "  -2.68%  "
$ws.Range("D14").Value = "'0.0000253"
$ws.Range("E14").Value = "  -3.70%  "
$ws.Range("D15").Value = "4.384.04"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").Value = "3.756.14"
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("D17").Value = "69.548.59"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("E18").Value = "  -1.55%  "
$ws.Range("E19").Value = "  -3.35%  "
$ws.Range("D20").Value = "'510.88"
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("D21").Value = "'16.41"
$ws.Range("E21").Value = "  -2.80%  "
$ws.Range("D22").Value = "'9.38"
$ws.Range("E22").Value = "  -1.41%  "
$ws.Range("D23").Value = "'0.727"
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("E24").Value = "  -0.74%  "
$ws.Range("D25").Value = "'86.49"
$ws.Range("E25").Value = "  -1.67%  "
$ws.Range("D26").Value = "'12.81"
$ws.Range("E26").Value = "  -3.27%  "
$ws.Range("D27").Value = "'0.0000137"
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("D28").Value = "'10.60"
$ws.Range("E28").Value = "  -4.70%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").Value = "'2.97"
$ws.Range("E31").Value = "  +3.60%  "
$ws.Range("D32").Value = "'7.94"
$ws.Range("E32").Value = "  +2.16%  "
$ws.Range("D33").Value = "'30.68"
$ws.Range("E33").Value = "  -2.50%  "
$ws.Range("E34").Value = "  -1.63%  "
$ws.Range("E36").Value = "  -0.61%  "
$ws.Range("E37").Value = "  -1.49%  "
$ws.Range("D38").Value = "'0.138"
$ws.Range("E38").Value = "  +2.87%  "
$ws.Range("E39").Value = "  +1.76%  "
$ws.Range("D40").Value = "'448.65"
$ws.Range("E40").Value = "  +6.50%  "
$ws.Range("E41").Value = "  -3.52%  "
$ws.Range("D42").Value = "'49.86"
$ws.Range("E42").Value = "  -2.34%  "
$ws.Range("D43").Value = "'2.95"
$ws.Range("E43").Value = "  +4.39%  "
$ws.Range("D44").Value = "'44.43"
$ws.Range("E44").Value = "  -1.68%  "
$ws.Range("D45").Value = "'8.57"
$ws.Range("D46").Value = "2.944.65"
$ws.Range("E46").Value = "  -3.25%  "
$ws.Range("D47").Value = "'0.0359"
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("D48").Value = "'27.22"
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "'138.58"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("E51").Value = "  -1.52%  "
